$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.834.18'
$ws.Range('E2').Value = '  +1.43%  '
$ws.Range('D3').Value = '2.120.62'
$ws.Range('E3').Value = '  +10.65%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '257.00'
$ws.Range('E5').Value = '  +3.18%  '
$ws.Range('D6').Value = '0.671'
$ws.Range('E6').Value = '  -3.23%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = '46.66'
$ws.Range('E8').Value = '  +6.68%  '
$ws.Range('D9').Value = '62.31'
$ws.Range('E9').Value = '  +8.44%  '
$ws.Range('D10').Value = '0.375'
$ws.Range('E10').Value = '  +2.80%  '
$ws.Range('D11').Value = '0.0745'
$ws.Range('E11').Value = '  -2.43%  '
$ws.Range('E12').Value = '  +0.83%  '
$ws.Range('D13').Value = '2.427.31'
$ws.Range('E13').Value = '  +10.37%  '
$ws.Range('D14').Value = '14.61'
$ws.Range('E14').Value = '  -0.06%  '
$ws.Range('D15').Value = '0.856'
$ws.Range('E15').Value = '  +7.47%  '
$ws.Range('D16').Value = '2.120.14'
$ws.Range('E16').Value = '  +10.64%  '
$ws.Range('E17').Value = '  +1.54%  '
$ws.Range('D18').Value = '36.836.28'
$ws.Range('E18').Value = '  +1.35%  '
$ws.Range('E19').Value = '  +0.23%  '
$ws.Range('D20').Value = '0.0₃0849'
$ws.Range('E20').Value = '  +0.47%  '
$ws.Range('D21').Value = '13.51'
$ws.Range('E21').Value = '  +2.39%  '
$ws.Range('D22').Value = '242.69'
$ws.Range('E22').Value = '  -4.06%  '
$ws.Range('D23').Value = '5.27'
$ws.Range('E23').Value = '  +1.32%  '
$ws.Range('E24').Value = '  +0.10%  '
$ws.Range('D25').Value = '2.50'
$ws.Range('E25').Value = '  -7.31%  '
$ws.Range('D26').Value = '173.10'
$ws.Range('E26').Value = '  +3.28%  '
$ws.Range('D27').Value = '21.47'
$ws.Range('E27').Value = '  +13.88%  '
$ws.Range('D28').Value = '9.26'
$ws.Range('E28').Value = '  +5.21%  '
$ws.Range('E29').Value = '  -8.90%  '
$ws.Range('E30').Value = '  -3.70%  '
$ws.Range('D31').Value = '23.00'
$ws.Range('E31').Value = '  +52.85%  '
$ws.Range('E32').Value = '  +0.84%  '
$ws.Range('D33').Value = '0.0965'
$ws.Range('E33').Value = '  +14.98%  '
$ws.Range('E34').Value = '  -0.58%  '
$ws.Range('D35').Value = '2.43'
$ws.Range('E35').Value = '  +20.93%  '
$ws.Range('D36').Value = '1.89'
$ws.Range('E36').Value = '  -3.27%  '
$ws.Range('E37').Value = '  -0.11%  '
$ws.Range('E38').Value = '  -2.73%  '
$ws.Range('D39').Value = '0.922'
$ws.Range('E39').Value = '  +7.16%  '
$ws.Range('D40').Value = '1.37'
$ws.Range('E40').Value = '  -7.71%  '
$ws.Range('D41').Value = '1.21'
$ws.Range('E41').Value = '  +8.75%  '
$ws.Range('E42').Value = '  -1.91%  '
$ws.Range('D43').Value = '99.55'
$ws.Range('E43').Value = '  -4.29%  '
$ws.Range('E44').Value = '  +17.51%  '
$ws.Range('E45').Value = '  -4.83%  '
$ws.Range('D46').Value = '1.368.03'
$ws.Range('E46').Value = '  +1.94%  '
$ws.Range('D47').Value = '0.0838'
$ws.Range('E47').Value = '  +3.72%  '
$ws.Range('D48').Value = '2.313.46'
$ws.Range('E48').Value = '  +10.18%  '
$ws.Range('E49').Value = '  -2.73%  '
$ws.Range('D50').Value = '6.92'
$ws.Range('E50').Value = '  +7.70%  '
$ws.Range('E51').Value = '  +1.97%  '
